$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: production-time values
$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 71
$ws.Range("C4").Value = 63
$ws.Range("D4").Value = 59

# Row 2: makespan label
$ws.Range("A2").Value = "make span: 150"

# Row 6: machine 1 orders
$ws.Range("A6").Value = "Order 9 - 48"
$ws.Range("B6").Value = "Order 5 - 40"
$ws.Range("C6").Value = "Order 1 - 61"
$ws.Range("D6").Value = "Order 8 - 34"

# Row 7: machine 2 orders
$ws.Range("A7").Value = "Order 3 - 74"
$ws.Range("B7").Value = "Order 2 - 79"
$ws.Range("C7").Value = "Order 6 - 87"
$ws.Range("D7").Value = "Order 7 - 91"

# Row 8: machine 3 orders - only column A now, C8 cleared
$ws.Range("A8").Value = "Order 4 - 80"
$ws.Range("C8").ClearContents()

# Row 9: new row with remaining order
$ws.Range("A9").Value = "Order 10 - 150"
